$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (edit specific substrings within rich-text shared strings) ---
$ws.Range("A8").Characters(21, 2).Text = "11"
$ws.Range("C9").Characters(27, 8).Text = "3/10/2025"
$ws.Range("C9").Characters(47, 8).Text = "3/16/2025"

# --- Cells that switch between numeric and text ("0"/"***.*") representations: copy from a donor cell that already has the desired style + shared-string value ---
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("F14").Copy($ws.Range("C16"))
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("C14").Copy($ws.Range("D19"))
$ws.Range("E14").Copy($ws.Range("E19"))
$ws.Range("J14").Copy($ws.Range("C20"))
$ws.Range("J14").Copy($ws.Range("F20"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("C14").Copy($ws.Range("C25"))
$ws.Range("C14").Copy($ws.Range("D25"))
$ws.Range("E14").Copy($ws.Range("E25"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))

# --- Plain numeric value updates ---
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 75
$ws.Range("L15").Value = 133.333333333333
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -45.454545454545
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 28
$ws.Range("K16").Value = -53.571428571428
$ws.Range("L16").Value = -45.833333333333
$ws.Range("M16").Value = -59.375
$ws.Range("N16").Value = -88.073394495412
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 42.857142857142
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -9.677419354838
$ws.Range("I17").Value = 55
$ws.Range("J17").Value = 66
$ws.Range("K17").Value = -16.666666666666
$ws.Range("L17").Value = 52.777777777777
$ws.Range("M17").Value = 66.666666666666
$ws.Range("N17").Value = -32.098765432098
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = -16.666666666666
$ws.Range("J18").Value = 17
$ws.Range("K18").Value = -47.058823529411
$ws.Range("L18").Value = -50
$ws.Range("N18").Value = -91.089108910891
$ws.Range("G19").Value = 12
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 28
$ws.Range("K19").Value = -12.5
$ws.Range("L19").Value = -17.647058823529
$ws.Range("M19").Value = 33.333333333333
$ws.Range("N19").Value = -46.153846153846
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("H20").Value = -83.333333333333
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = -69.230769230769
$ws.Range("L20").Value = -75
$ws.Range("M20").Value = -78.947368421052
$ws.Range("N20").Value = -95.876288659793
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 13.333333333333
$ws.Range("F21").Value = 57
$ws.Range("H21").Value = -16.176470588235
$ws.Range("I21").Value = 118
$ws.Range("J21").Value = 161
$ws.Range("K21").Value = -26.708074534161
$ws.Range("L21").Value = -9.923664122137
$ws.Range("M21").Value = -11.940298507462
$ws.Range("N21").Value = -73.719376391982
$ws.Range("L22").Value = -80
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -28.571428571428
$ws.Range("I23").Value = 12
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = -29.411764705882
$ws.Range("L23").Value = -29.411764705882
$ws.Range("M23").Value = 0
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 56
$ws.Range("G24").Value = 47
$ws.Range("H24").Value = 19.148936170212
$ws.Range("I24").Value = 138
$ws.Range("J24").Value = 129
$ws.Range("K24").Value = 6.976744186046
$ws.Range("L24").Value = 45.263157894736
$ws.Range("M24").Value = 115.625
$ws.Range("F25").Value = 9
$ws.Range("H25").Value = 200
$ws.Range("L25").Value = 30.769230769230
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 12.5
$ws.Range("F26").Value = 36
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = 71.428571428571
$ws.Range("I26").Value = 94
$ws.Range("J26").Value = 76
$ws.Range("K26").Value = 23.684210526315
$ws.Range("L26").Value = 44.615384615384
$ws.Range("M26").Value = 22.077922077922
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = 60
$ws.Range("L27").Value = 14.285714285714
